$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, shifting existing rows 49..161 down to 50..162
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record
$ws.Range("A49").Value = 11
$ws.Range("B49").Value = "Vega Monumental Concepción"
$ws.Range("C49").Value = "Bíobío"
$ws.Range("D49").Value = 45281
$ws.Range("E49").Value = 8
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100101
$ws.Range("H49").Value = "Berries"
$ws.Range("I49").Value = 100101001
$ws.Range("J49").Value = "Arándano (blue)"
$ws.Range("K49").Value = "Sin especificar"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 180
$ws.Range("N49").Value = 4000
$ws.Range("O49").Value = 4500
$ws.Range("P49").Value = 4278
$ws.Range("Q49").Value = "`$/bandeja 2 kilos"
$ws.Range("R49").Value = "Región de Ñuble"
$ws.Range("S49").Value = 2139
$ws.Range("T49").Value = 2
